$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta")

# The "meta" sheet stores key/value metadata pairs in column A/B, one pair
# per row, with column A styled (s="1", bold orange). Row 8 previously held
# just a trailing, empty, styled cell (A8) with nothing in B8 (used as an
# end marker). We now add a new "style" / "default" key-value pair in row 8,
# and push the old trailing empty styled cell down to row 9.

# First, move the formatting of the current trailing empty cell (A8) down to
# the new trailing row (A9), preserving its exact style (no value).
$ws.Range("A8").Copy($ws.Range("A9"))

# Now populate row 8 with the new "style"/"default" metadata pair.
$ws.Range("A8").Value = "style"
$ws.Range("B8").Value = "default"
